$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> (new D value or $null, new E value)
$changes = @(
    @{ Row = 2;  D = "27.428.84";     E = "  +1.89%  " },
    @{ Row = 3;  D = "1.842.18";      E = "  +1.43%  " },
    @{ Row = 4;  D = $null;           E = "  +1.40%  " },
    @{ Row = 5;  D = "315.16";        E = $null },
    @{ Row = 6;  D = $null;           E = "  +1.25%  " },
    @{ Row = 7;  D = "0.4767";        E = "  +1.69%  " },
    @{ Row = 8;  D = $null;           E = "  +0.42%  " },
    @{ Row = 9;  D = "0.07471";       E = "  +1.20%  " },
    @{ Row = 10; D = "0.8871";        E = "  +1.71%  " },
    @{ Row = 11; D = "20.50";         E = "  +0.48%  " },
    @{ Row = 12; D = "1.854.40";      E = "  +1.90%  " },
    @{ Row = 13; D = "0.07369";       E = "  +4.10%  " },
    @{ Row = 14; D = "5.488";         E = "  +1.85%  " },
    @{ Row = 15; D = "93.27";         E = "  +1.59%  " },
    @{ Row = 16; D = "6.601";         E = "  +1.12%  " },
    @{ Row = 17; D = "1.017";         E = "  +1.36%  " },
    @{ Row = 18; D = "0.000008859";   E = "  +1.51%  " },
    @{ Row = 19; D = $null;           E = "  +1.34%  " },
    @{ Row = 20; D = "14.84";         E = "  +0.63%  " },
    @{ Row = 21; D = "27.442.41";     E = "  +1.79%  " },
    @{ Row = 22; D = "5.349";         E = "  +0.42%  " },
    @{ Row = 23; D = $null;           E = "  +1.20%  " },
    @{ Row = 24; D = "2.084.51";      E = "  +1.38%  " },
    @{ Row = 25; D = "1.910";         E = "  +1.03%  " },
    @{ Row = 26; D = "152.58";        E = "  +1.15%  " },
    @{ Row = 27; D = $null;           E = "  +1.61%  " },
    @{ Row = 28; D = "2.170";         E = "  -0.22%  " },
    @{ Row = 29; D = "5.287";         E = "  -1.03%  " },
    @{ Row = 30; D = "118.20";        E = "  +1.74%  " },
    @{ Row = 31; D = "0.08987";       E = "  +0.23%  " },
    @{ Row = 32; D = "0.7592";        E = "  -1.39%  " },
    @{ Row = 33; D = $null;           E = "  +1.20%  " },
    @{ Row = 34; D = "4.572";         E = "  +1.26%  " },
    @{ Row = 35; D = "2.952";         E = "  +1.40%  " },
    @{ Row = 36; D = $null;           E = "  +1.37%  " },
    @{ Row = 37; D = "1.107";         E = "  +1.96%  " },
    @{ Row = 38; D = "0.05365";       E = "  +1.14%  " },
    @{ Row = 39; D = "0.01965";       E = "  +0.04%  " },
    @{ Row = 40; D = "3.002";         E = "  +1.99%  " },
    @{ Row = 41; D = "7.320";         E = "  +0.62%  " },
    @{ Row = 42; D = "0.5363";        E = "  +0.34%  " },
    @{ Row = 43; D = "2.378";         E = "  +1.62%  " },
    @{ Row = 44; D = "0.1667";        E = "  +0.55%  " },
    @{ Row = 45; D = "8.558";         E = "  +1.18%  " },
    @{ Row = 46; D = "0.4978";        E = "  +0.89%  " },
    @{ Row = 47; D = $null;           E = "  +1.34%  " },
    @{ Row = 48; D = $null;           E = "  +1.40%  " },
    @{ Row = 49; D = "105.37";        E = "  +1.76%  " },
    @{ Row = 50; D = $null;           E = "  +0.73%  " },
    @{ Row = 51; D = "0.06328";       E = "  +0.38%  " }
)

foreach ($change in $changes) {
    $r = $change.Row
    if ($null -ne $change.D) {
        # Leading apostrophe forces Excel to store this as text rather than
        # auto-coercing a numeric-looking string (e.g. "315.16") into a
        # floating point number.
        $ws.Range("D$r").Value = "'" + $change.D
    }
    if ($null -ne $change.E) {
        $ws.Range("E$r").Value = $change.E
    }
}
